$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.882.37"

$ws.Range("D3").Value = "1.634.67"

$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.57"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").Value = "  -0.43%  "

$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("D10").ClearFormats()

$ws.Range("D12").Value = "1.859.43"

$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").Value = "1.635.95"

$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.56"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "25.895.17"

$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.69"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.93"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.56"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("E27").Value = "  +3.09%  "

$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.42"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("E35").Value = "  +1.46%  "

$ws.Range("D37").Value = "1.139.15"

$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("E40").Value = "  +0.58%  "

$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.43"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D45").Value = "1.768.68"

$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("E46").Value = "  +1.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.31"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "  +1.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0529"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = "  +3.31%  "

$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("E50").Value = "  +1.55%  "

$ws.Range("E51").Value = "  -0.19%  "
